$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 7374.375
$ws.Range("I76").Value = 5000
$ws.Range("J76").Value = 7713.5713
$ws.Range("K76").Value = 5000
$ws.Range("L76").Value = 7713.5713
$ws.Range("M76").Value = -4685
$ws.Range("N76").Value = -8343.5713
$ws.Range("H79").Value = 7374.375
$ws.Range("I79").Value = 5000
$ws.Range("J79").Value = 7713.5713
$ws.Range("K79").Value = 5000
$ws.Range("L79").Value = 7713.5713
$ws.Range("M79").Value = -3908
$ws.Range("N79").Value = -9897.5713
$ws.Range("H86").Value = 6257.0713
$ws.Range("I86").Value = 1516.5
$ws.Range("K86").Value = 1516.5
$ws.Range("M86").Value = -393.5
$ws.Range("H88").Value = 5499.75
$ws.Range("J88").Value = 4000
$ws.Range("L88").Value = 4000
$ws.Range("N88").Value = -4812
$ws.Range("H89").Value = 6257.0713
$ws.Range("I89").Value = 1516.5
$ws.Range("K89").Value = 7582.5
$ws.Range("M89").Value = -1966.5
$ws.Range("H91").Value = 5499.75
$ws.Range("J91").Value = 4000
$ws.Range("L91").Value = 4000
$ws.Range("M91").Value = -4595.6665
$ws.Range("N91").Value = -6808
$ws.Range("H99").Value = 2601.6
$ws.Range("I99").Value = 2788.889
$ws.Range("J99").Value = 2320.6667
$ws.Range("K99").Value = 8366.667000000001
$ws.Range("L99").Value = 6962.000100000001
$ws.Range("M99").Value = -6868.667000000001
$ws.Range("N99").Value = -9958.000100000001
$ws.Range("H101").Value = 831.13336
$ws.Range("I101").Value = 875.5454999999999
$ws.Range("J101").Value = 709
$ws.Range("K101").Value = 2626.6365
$ws.Range("L101").Value = 2127
$ws.Range("M101").Value = -1004.6365
$ws.Range("N101").Value = -5371
$ws.Range("H115").Value = 640.75
$ws.Range("I115").Value = 640.75
$ws.Range("K115").Value = 1922.25
$ws.Range("M115").Value = -355.25
$ws.Range("H118").Value = 1217
$ws.Range("I118").Value = 1073.5714
$ws.Range("J118").Value = 2221
$ws.Range("K118").Value = 3220.7142
$ws.Range("L118").Value = 6663
$ws.Range("M118").Value = -1563.7142
$ws.Range("N118").Value = -9977
$ws.Range("H127").Value = 2973.25
$ws.Range("I127").Value = 2973.25
$ws.Range("K127").Value = 8919.75
$ws.Range("M127").Value = -3959.75
$ws.Range("H129").Value = 18779.334
$ws.Range("I129").Value = 12535.2
$ws.Range("J129").Value = 50000
$ws.Range("K129").Value = 37605.60000000001
$ws.Range("L129").Value = 150000
$ws.Range("M129").Value = -32605.60000000001
$ws.Range("N129").Value = -160000
$ws.Range("H138").Value = 4279.7334
$ws.Range("J138").Value = 6308.6313
$ws.Range("L138").Value = 18925.8939
$ws.Range("N138").Value = -29205.8939

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4759.9614
$ws.Range("I32").Value = 1853.3972
$ws.Range("K32").Value = 1853.3972
$ws.Range("M32").Value = -1566.3972
$ws.Range("H45").Value = 1974.7727
$ws.Range("I45").Value = 1444.6842
$ws.Range("J45").Value = 5332
$ws.Range("K45").Value = 1444.6842
$ws.Range("L45").Value = 5332
$ws.Range("M45").Value = -1067.6842
$ws.Range("N45").Value = -6086
$ws.Range("H61").Value = 8989.65
$ws.Range("I61").Value = 7730.758
$ws.Range("K61").Value = 7730.758
$ws.Range("M61").Value = -7518.758
$ws.Range("H132").Value = 3540.1292
$ws.Range("I132").Value = 1773.762
$ws.Range("J132").Value = 7249.5
$ws.Range("K132").Value = 5321.286
$ws.Range("L132").Value = 21748.5
$ws.Range("M132").Value = -2791.286
$ws.Range("N132").Value = -26808.5
$ws.Range("H136").Value = 8989.65
$ws.Range("I136").Value = 7730.758
$ws.Range("K136").Value = 23192.274
$ws.Range("M136").Value = -20642.274

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 850.4
$ws.Range("I94").Value = 561.86365
$ws.Range("K94").Value = 561.86365
$ws.Range("M94").Value = -110.86365
$ws.Range("H134").Value = 5181.8184
$ws.Range("I134").Value = 2000
$ws.Range("K134").Value = 6000
$ws.Range("M134").Value = -3465

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1806.5
$ws.Range("I122").Value = 1817.3636
$ws.Range("K122").Value = 5452.0908
$ws.Range("M122").Value = -3002.0908

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 510.04544
$ws.Range("I114").Value = 383.45456
$ws.Range("K114").Value = 1150.36368
$ws.Range("M114").Value = 2103.63632
$ws.Range("H117").Value = 1315.5
$ws.Range("I117").Value = 720.8
$ws.Range("K117").Value = 2162.4
$ws.Range("M117").Value = 1279.6
$ws.Range("H131").Value = 588189.5
$ws.Range("I131").Value = 914.8889
$ws.Range("K131").Value = 2744.6667
$ws.Range("M131").Value = 2295.3333
$ws.Range("H137").Value = 3992.7083
$ws.Range("I137").Value = 2253.4546
$ws.Range("K137").Value = 6760.3638
$ws.Range("M137").Value = -1660.3638
$ws.Range("H140").Value = 1628.4
$ws.Range("I140").Value = 1443.8889
$ws.Range("J140").Value = 1692.2693
$ws.Range("K140").Value = 4331.6667
$ws.Range("L140").Value = 5076.8079
$ws.Range("M140").Value = 848.3333000000002
$ws.Range("N140").Value = -15436.8079

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5021.891
$ws.Range("I70").Value = 5350
$ws.Range("K70").Value = 5350
$ws.Range("M70").Value = -5080
$ws.Range("H73").Value = 5021.891
$ws.Range("I73").Value = 5350
$ws.Range("K73").Value = 5350
$ws.Range("M73").Value = -4414
$ws.Range("H97").Value = 607.125
$ws.Range("I97").Value = 736.1111
$ws.Range("J97").Value = 441.2857
$ws.Range("K97").Value = 736.1111
$ws.Range("L97").Value = 441.2857
$ws.Range("M97").Value = -240.1111
$ws.Range("N97").Value = -1433.2857
$ws.Range("H122").Value = 5557.2915
$ws.Range("I122").Value = 4234.9
$ws.Range("K122").Value = 12704.7
$ws.Range("M122").Value = -10254.7
$ws.Range("H132").Value = 5958.6177
$ws.Range("I132").Value = 5912.273
$ws.Range("K132").Value = 17736.819
$ws.Range("M132").Value = -15206.819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4282.636
$ws.Range("J46").Value = 4789.8887
$ws.Range("L46").Value = 4789.8887
$ws.Range("N46").Value = -5165.8887
$ws.Range("H68").Value = 8672.450000000001
$ws.Range("I68").Value = 8229.933999999999
$ws.Range("J68").Value = 10000
$ws.Range("K68").Value = 8229.933999999999
$ws.Range("L68").Value = 10000
$ws.Range("M68").Value = -7480.933999999999
$ws.Range("N68").Value = -11498
$ws.Range("H71").Value = 8672.450000000001
$ws.Range("I71").Value = 8229.933999999999
$ws.Range("J71").Value = 10000
$ws.Range("K71").Value = 41149.67
$ws.Range("L71").Value = 50000
$ws.Range("M71").Value = -37405.67
$ws.Range("N71").Value = -57488
$ws.Range("H93").Value = 1799.7084
$ws.Range("I93").Value = 1687.9412
$ws.Range("J93").Value = 2071.1428
$ws.Range("K93").Value = 1687.9412
$ws.Range("L93").Value = 2071.1428
$ws.Range("M93").Value = -439.9412
$ws.Range("N93").Value = -4567.1428

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3200.3635
$ws.Range("I96").Value = 1037.5
$ws.Range("J96").Value = 5795.8
$ws.Range("K96").Value = 1037.5
$ws.Range("L96").Value = 5795.8
$ws.Range("M96").Value = 335.5
$ws.Range("N96").Value = -8541.799999999999
$ws.Range("H132").Value = 3094.4783
$ws.Range("I132").Value = 1342.0476
$ws.Range("J132").Value = 21495
$ws.Range("K132").Value = 4026.142800000001
$ws.Range("L132").Value = 64485
$ws.Range("M132").Value = -1496.142800000001
$ws.Range("N132").Value = -69545
